# Automatic update of files.
# Update the "Förändrad" (changed) date in column C for rows 2-14
# from 45174 (2023-09-05) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C14").Value = 45175
